$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.225934386253357
$ws.Range("B1").Value = 1.68542218208313
$ws.Range("C1").Value = 2.969988107681274
$ws.Range("D1").Value = 1.501512765884399
$ws.Range("E1").Value = 0.8209075331687927
